# Commit: Modif url canonique termino df9498eb894642b7264f6d5c9a38a249f1b02b34
$wb = $excel.ActiveWorkbook

# --- 1) Metadata sheet: bump the IG "Date" property ---
$meta = $wb.Worksheets.Item("Metadata")
$metaUsed = $meta.UsedRange
$metaUsed.Replace("2025-07-24T13:17:05+00:00", "2025-07-25T07:22:51+00:00") | Out-Null

# --- 2) Elements sheet: rewrite the canonical terminology (CodeSystem) URLs ---
#     https://interop.esante.gouv.fr/terminologies/CodeSystem-<ID>?vs
#  -> https://mos.esante.gouv.fr/NOS/<ID-with-first-dash-as-underscore>/FHIR/<ID>?vs
$els = $wb.Worksheets.Item("Elements")
$elsUsed = $els.UsedRange

$elsUsed.Replace("https://interop.esante.gouv.fr/terminologies/CodeSystem-TRE-R35-TypeVoie?vs", "https://mos.esante.gouv.fr/NOS/TRE_R35-TypeVoie/FHIR/TRE-R35-TypeVoie?vs") | Out-Null
$elsUsed.Replace("https://interop.esante.gouv.fr/terminologies/CodeSystem-TRE-R13-CommuneOM?vs", "https://mos.esante.gouv.fr/NOS/TRE_R13-CommuneOM/FHIR/TRE-R13-CommuneOM?vs") | Out-Null
$elsUsed.Replace("https://interop.esante.gouv.fr/terminologies/CodeSystem-TRE-R14-TypeDiplome?vs", "https://mos.esante.gouv.fr/NOS/TRE_R14-TypeDiplome/FHIR/TRE-R14-TypeDiplome?vs") | Out-Null
$elsUsed.Replace("https://interop.esante.gouv.fr/terminologies/CodeSystem-TRE-R16-LieuFormation?vs", "https://mos.esante.gouv.fr/NOS/TRE_R16-LieuFormation/FHIR/TRE-R16-LieuFormation?vs") | Out-Null
$elsUsed.Replace("https://interop.esante.gouv.fr/terminologies/CodeSystem-TRE-R56-Attestation?vs", "https://mos.esante.gouv.fr/NOS/TRE_R56-Attestation/FHIR/TRE-R56-Attestation?vs") | Out-Null
$elsUsed.Replace("https://interop.esante.gouv.fr/terminologies/CodeSystem-TRE-R36-AutreDiplomeObtenu?vs", "https://mos.esante.gouv.fr/NOS/TRE_R36-AutreDiplomeObtenu/FHIR/TRE-R36-AutreDiplomeObtenu?vs") | Out-Null
$elsUsed.Replace("https://interop.esante.gouv.fr/terminologies/CodeSystem-TRE-R257-TypeBAL?vs", "https://mos.esante.gouv.fr/NOS/TRE_R257-TypeBAL/FHIR/TRE-R257-TypeBAL?vs") | Out-Null
$elsUsed.Replace("https://interop.esante.gouv.fr/terminologies/CodeSystem-TRE-R52-CapaciteDiplome?vs", "https://mos.esante.gouv.fr/NOS/TRE_R52-CapaciteDiplome/FHIR/TRE-R52-CapaciteDiplome?vs") | Out-Null
$elsUsed.Replace("https://interop.esante.gouv.fr/terminologies/CodeSystem-TRE-R04-TypeSavoirFaire?vs", "https://mos.esante.gouv.fr/NOS/TRE_R04-TypeSavoirFaire/FHIR/TRE-R04-TypeSavoirFaire?vs") | Out-Null
$elsUsed.Replace("https://interop.esante.gouv.fr/terminologies/CodeSystem-TRE-R55-CertificatEtudeSpeciale?vs", "https://mos.esante.gouv.fr/NOS/TRE_R55-CertificatEtudeSpeciale/FHIR/TRE-R55-CertificatEtudeSpeciale?vs") | Out-Null
$elsUsed.Replace("https://interop.esante.gouv.fr/terminologies/CodeSystem-TRE-R50-DESCGroupe1Diplome?vs", "https://mos.esante.gouv.fr/NOS/TRE_R50-DESCGroupe1Diplome/FHIR/TRE-R50-DESCGroupe1Diplome?vs") | Out-Null
$elsUsed.Replace("https://interop.esante.gouv.fr/terminologies/CodeSystem-TRE-R51-DESCGroupe2Diplome?vs", "https://mos.esante.gouv.fr/NOS/TRE_R51-DESCGroupe2Diplome/FHIR/TRE-R51-DESCGroupe2Diplome?vs") | Out-Null
$elsUsed.Replace("https://interop.esante.gouv.fr/terminologies/CodeSystem-TRE-R226-Dip2iemeCycleNQ?vs", "https://mos.esante.gouv.fr/NOS/TRE_R226-Dip2iemeCycleNQ/FHIR/TRE-R226-Dip2iemeCycleNQ?vs") | Out-Null
$elsUsed.Replace("https://interop.esante.gouv.fr/terminologies/CodeSystem-TRE-R58-AutreTypeDiplome?vs", "https://mos.esante.gouv.fr/NOS/TRE_R58-AutreTypeDiplome/FHIR/TRE-R58-AutreTypeDiplome?vs") | Out-Null
$elsUsed.Replace("https://interop.esante.gouv.fr/terminologies/CodeSystem-TRE-R53-DiplomePaysEEE?vs", "https://mos.esante.gouv.fr/NOS/TRE_R53-DiplomePaysEEE/FHIR/TRE-R53-DiplomePaysEEE?vs") | Out-Null
$elsUsed.Replace("https://interop.esante.gouv.fr/terminologies/CodeSystem-TRE-R57-DiplomeEuropeenEtudeSpecialisee?vs", "https://mos.esante.gouv.fr/NOS/TRE_R57-DiplomeEuropeenEtudeSpecialisee/FHIR/TRE-R57-DiplomeEuropeenEtudeSpecialisee?vs") | Out-Null
$elsUsed.Replace("https://interop.esante.gouv.fr/terminologies/CodeSystem-TRE-R48-DiplomeEtatFrancais?vs", "https://mos.esante.gouv.fr/NOS/TRE_R48-DiplomeEtatFrancais/FHIR/TRE-R48-DiplomeEtatFrancais?vs") | Out-Null
$elsUsed.Replace("https://interop.esante.gouv.fr/terminologies/CodeSystem-TRE-R49-DiplomeEtudeSpecialisee?vs", "https://mos.esante.gouv.fr/NOS/TRE_R49-DiplomeEtudeSpecialisee/FHIR/TRE-R49-DiplomeEtudeSpecialisee?vs") | Out-Null
$elsUsed.Replace("https://interop.esante.gouv.fr/terminologies/CodeSystem-TRE-R54-DiplomeUniversiteInterUniversitaire?vs", "https://mos.esante.gouv.fr/NOS/TRE_R54-DiplomeUniversiteInterUniversitaire/FHIR/TRE-R54-DiplomeUniversiteInterUniversitaire?vs") | Out-Null
$elsUsed.Replace("https://interop.esante.gouv.fr/terminologies/CodeSystem-TRE-R11-CiviliteExercice?vs", "https://mos.esante.gouv.fr/NOS/TRE_R11-CiviliteExercice/FHIR/TRE-R11-CiviliteExercice?vs") | Out-Null
$elsUsed.Replace("https://interop.esante.gouv.fr/terminologies/CodeSystem-TRE-R09-CategorieProfessionnelle?vs", "https://mos.esante.gouv.fr/NOS/TRE_R09-CategorieProfessionnelle/FHIR/TRE-R09-CategorieProfessionnelle?vs") | Out-Null
$elsUsed.Replace("https://interop.esante.gouv.fr/terminologies/CodeSystem-TRE-G09-DepartementOM?vs", "https://mos.esante.gouv.fr/NOS/TRE_G09-DepartementOM/FHIR/TRE-G09-DepartementOM?vs") | Out-Null
$elsUsed.Replace("https://interop.esante.gouv.fr/terminologies/CodeSystem-TRE-R82-Ordre?vs", "https://mos.esante.gouv.fr/NOS/TRE_R82-Ordre/FHIR/TRE-R82-Ordre?vs") | Out-Null
$elsUsed.Replace("https://interop.esante.gouv.fr/terminologies/CodeSystem-TRE-R33-StatutInscription?vs", "https://mos.esante.gouv.fr/NOS/TRE_R33-StatutInscription/FHIR/TRE-R33-StatutInscription?vs") | Out-Null
$elsUsed.Replace("https://interop.esante.gouv.fr/terminologies/CodeSystem-TRE-R03-AttributionParticuliere?vs", "https://mos.esante.gouv.fr/NOS/TRE_R03-AttributionParticuliere/FHIR/TRE-R03-AttributionParticuliere?vs") | Out-Null
$elsUsed.Replace("https://interop.esante.gouv.fr/terminologies/CodeSystem-TRE-R223-NatCycleForm?vs", "https://mos.esante.gouv.fr/NOS/TRE_R223-NatCycleForm/FHIR/TRE-R223-NatCycleForm?vs") | Out-Null
$elsUsed.Replace("https://interop.esante.gouv.fr/terminologies/CodeSystem-TRE-R224-NiveauFormAcquis?vs", "https://mos.esante.gouv.fr/NOS/TRE_R224-NiveauFormAcquis/FHIR/TRE-R224-NiveauFormAcquis?vs") | Out-Null
$elsUsed.Replace("https://interop.esante.gouv.fr/terminologies/CodeSystem-TRE-R225-AnneeUniversitaire?vs", "https://mos.esante.gouv.fr/NOS/TRE_R225-AnneeUniversitaire/FHIR/TRE-R225-AnneeUniversitaire?vs") | Out-Null
$elsUsed.Replace("https://interop.esante.gouv.fr/terminologies/CodeSystem-TRE-R348-FormationSpecialiseeTransversale?vs", "https://mos.esante.gouv.fr/NOS/TRE_R348-FormationSpecialiseeTransversale/FHIR/TRE-R348-FormationSpecialiseeTransversale?vs") | Out-Null
$elsUsed.Replace("https://interop.esante.gouv.fr/terminologies/CodeSystem-TRE-R252-TypeHoraire?vs", "https://mos.esante.gouv.fr/NOS/TRE_R252-TypeHoraire/FHIR/TRE-R252-TypeHoraire?vs") | Out-Null
$elsUsed.Replace("https://interop.esante.gouv.fr/terminologies/CodeSystem-TRE-G00-Langue?vs", "https://mos.esante.gouv.fr/NOS/TRE_G00-Langue/FHIR/TRE-G00-Langue?vs") | Out-Null
$elsUsed.Replace("https://interop.esante.gouv.fr/terminologies/CodeSystem-TRE-R81-Civilite?vs", "https://mos.esante.gouv.fr/NOS/TRE_R81-Civilite/FHIR/TRE-R81-Civilite?vs") | Out-Null
$elsUsed.Replace("https://interop.esante.gouv.fr/terminologies/CodeSystem-TRE-R10-SexeAdministratif?vs", "https://mos.esante.gouv.fr/NOS/TRE_R10-SexeAdministratif/FHIR/TRE-R10-SexeAdministratif?vs") | Out-Null
$elsUsed.Replace("https://interop.esante.gouv.fr/terminologies/CodeSystem-TRE-R31-StatutEtatCivil?vs", "https://mos.esante.gouv.fr/NOS/TRE_R31-StatutEtatCivil/FHIR/TRE-R31-StatutEtatCivil?vs") | Out-Null
$elsUsed.Replace("https://interop.esante.gouv.fr/terminologies/CodeSystem-TRE-R20-Pays?vs", "https://mos.esante.gouv.fr/NOS/TRE_R20-Pays/FHIR/TRE-R20-Pays?vs") | Out-Null
$elsUsed.Replace("https://interop.esante.gouv.fr/terminologies/CodeSystem-TRE-G08-TypeIdentifiantPersonne?vs", "https://mos.esante.gouv.fr/NOS/TRE_G08-TypeIdentifiantPersonne/FHIR/TRE-G08-TypeIdentifiantPersonne?vs") | Out-Null
$elsUsed.Replace("https://interop.esante.gouv.fr/terminologies/CodeSystem-TRE-R23-ModeExercice?vs", "https://mos.esante.gouv.fr/NOS/TRE_R23-ModeExercice/FHIR/TRE-R23-ModeExercice?vs") | Out-Null
$elsUsed.Replace("https://interop.esante.gouv.fr/terminologies/CodeSystem-TRE-R22-GenreActivite?vs", "https://mos.esante.gouv.fr/NOS/TRE_R22-GenreActivite/FHIR/TRE-R22-GenreActivite?vs") | Out-Null
$elsUsed.Replace("https://interop.esante.gouv.fr/terminologies/CodeSystem-TRE-R25-MotifFinActivite?vs", "https://mos.esante.gouv.fr/NOS/TRE_R25-MotifFinActivite/FHIR/TRE-R25-MotifFinActivite?vs") | Out-Null
$elsUsed.Replace("https://interop.esante.gouv.fr/terminologies/CodeSystem-TRE-R32-StatutHospitalier?vs", "https://mos.esante.gouv.fr/NOS/TRE_R32-StatutHospitalier/FHIR/TRE-R32-StatutHospitalier?vs") | Out-Null
$elsUsed.Replace("https://interop.esante.gouv.fr/terminologies/CodeSystem-TRE-R06-SectionTableauCNOP?vs", "https://mos.esante.gouv.fr/NOS/TRE_R06-SectionTableauCNOP/FHIR/TRE-R06-SectionTableauCNOP?vs") | Out-Null
$elsUsed.Replace("https://interop.esante.gouv.fr/terminologies/CodeSystem-TRE-G05-SousSectionTableauCNOP?vs", "https://mos.esante.gouv.fr/NOS/TRE_G05-SousSectionTableauCNOP/FHIR/TRE-G05-SousSectionTableauCNOP?vs") | Out-Null
$elsUsed.Replace("https://interop.esante.gouv.fr/terminologies/CodeSystem-TRE-R24-TypeActiviteLiberale?vs", "https://mos.esante.gouv.fr/NOS/TRE_R24-TypeActiviteLiberale/FHIR/TRE-R24-TypeActiviteLiberale?vs") | Out-Null
$elsUsed.Replace("https://interop.esante.gouv.fr/terminologies/CodeSystem-TRE-R34-StatutProfessionnelSSA?vs", "https://mos.esante.gouv.fr/NOS/TRE_R34-StatutProfessionnelSSA/FHIR/TRE-R34-StatutProfessionnelSSA?vs") | Out-Null
$elsUsed.Replace("https://interop.esante.gouv.fr/terminologies/CodeSystem-TRE-R243-CompetenceSpecifique?vs", "https://mos.esante.gouv.fr/NOS/TRE_R243-CompetenceSpecifique/FHIR/TRE-R243-CompetenceSpecifique?vs") | Out-Null
$elsUsed.Replace("https://interop.esante.gouv.fr/terminologies/CodeSystem-TRE-R282-CNAMAmeliSecteurConventionnement?vs", "https://mos.esante.gouv.fr/NOS/TRE_R282-CNAMAmeliSecteurConventionnement/FHIR/TRE-R282-CNAMAmeliSecteurConventionnement?vs") | Out-Null
$elsUsed.Replace("https://interop.esante.gouv.fr/terminologies/CodeSystem-TRE-R200-CanalCommunication?vs", "https://mos.esante.gouv.fr/NOS/TRE_R200-CanalCommunication/FHIR/TRE-R200-CanalCommunication?vs") | Out-Null
$elsUsed.Replace("https://interop.esante.gouv.fr/terminologies/CodeSystem-TRE-R256-TypeMessagerie?vs", "https://mos.esante.gouv.fr/NOS/TRE_R256-TypeMessagerie/FHIR/TRE-R256-TypeMessagerie?vs") | Out-Null
$elsUsed.Replace("https://interop.esante.gouv.fr/terminologies/CodeSystem-TRE-R283-NiveauConfidentialite?vs", "https://mos.esante.gouv.fr/NOS/TRE_R283-NiveauConfidentialite/FHIR/TRE-R283-NiveauConfidentialite?vs") | Out-Null

# --- 3) Widen column AA ("Binding Value Set") to fit the longer URLs ---
$bindingCol = $els.Columns.Item(27)
$bindingCol.ColumnWidth = 109.8333

